$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2; D=45097; M=90; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=3; D=45099; M=200; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=4; D=44431; M=100; N=1300; O=1300; P=1300; S=1300 }
    @{ Row=5; D=45079; M=30; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=6; D=44830; M=50; N=2500; O=2500; P=2500; S=2500 }
    @{ Row=7; D=45075; M=240; N=3200; O=3200; P=3200; S=3200 }
    @{ Row=8; D=44418; M=40; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=9; D=45042; M=25; N=3500; O=3500; P=3500; S=3500 }
    @{ Row=10; D=44476; M=80; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=11; D=45092; M=120; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=12; D=44357; M=35; N=1000; O=1000; P=1000; S=1000 }
    @{ Row=13; D=45085; M=40; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=14; D=44762; M=50; N=2300; O=2300; P=2300; S=2300 }
    @{ Row=15; D=45093; M=90; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=16; D=44749; M=120; N=2300; O=2300; P=2300; S=2300 }
    @{ Row=17; D=45068; M=50; N=3250; O=3250; P=3250; S=3250 }
    @{ Row=18; D=44748; M=300; N=2300; O=2300; P=2300; S=2300 }
    @{ Row=19; D=44811; M=60; N=2500; O=2500; P=2500; S=2500 }
    @{ Row=20; D=44473; M=120; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=21; D=45104; M=50; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=22; D=45054; M=25; N=2500; O=2500; P=2500; S=2500 }
    @{ Row=23; D=44405; M=50; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=24; D=44812; M=50; N=2500; O=2500; P=2500; S=2500 }
    @{ Row=25; D=45055; M=25; N=2800; O=2800; P=2800; S=2800 }
    @{ Row=26; D=44417; M=80; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=27; D=44432; M=30; N=1300; O=1300; P=1300; S=1300 }
    @{ Row=28; D=45086; M=30; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=29; D=44424; M=50; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=30; D=44343; M=60; N=1300; O=1300; P=1300; S=1300 }
    @{ Row=31; D=45062; M=60; N=3200; O=3200; P=3200; S=3200 }
    @{ Row=32; D=45090; M=50; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=33; D=44753; M=160; N=2300; O=2300; P=2300; S=2300 }
    @{ Row=34; D=44438; M=60; N=1200; O=1200; P=1200; S=1200 }
    @{ Row=35; D=45112; M=50; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=36; D=45044; M=150; N=3500; O=3500; P=3500; S=3500 }
    @{ Row=37; D=45106; M=120; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=38; D=44435; M=130; N=1300; O=1300; P=1300; S=1300 }
    @{ Row=39; D=44763; M=50; N=2300; O=2300; P=2300; S=2300 }
    @{ Row=40; D=45111; M=50; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=41; D=45149; M=100; N=2700; O=2700; P=2700; S=2700 }
    @{ Row=42; D=44760; M=80; N=2300; O=2300; P=2300; S=2300 }
    @{ Row=43; D=45113; M=90; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=44; D=45076; M=100; N=2600; O=2600; P=2600; S=2600 }
    @{ Row=45; D=45148; M=280; N=2750; O=2750; P=2750; S=2750 }
    @{ Row=46; D=45041; M=80; N=3500; O=3500; P=3500; S=3500 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D   # Column D: Fecha
    $ws.Cells.Item($r, 13).Value = $item.M  # Column M: Volumen
    $ws.Cells.Item($r, 14).Value = $item.N  # Column N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O  # Column O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P  # Column P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $item.S  # Column S: Precio $/Kg
}
